$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TEST_CASES")
$ws2 = $wb.Worksheets.Item("STEPS")

# --- Fill in the previously-empty TC_REFERENCE (G) / TC_NAME (H) columns ---
# The parser used to skip these optional columns; now they are populated.
$refs  = @("ref1","ref2","ref3","ref4","ref5","ref6","ref7","ref8","ref9","ref10")
$names = @("name1","name2","name3","name4","name5","name6","name7","name8","name9","name10")

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 7).Value = $refs[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 8).Value = $names[$i]
}

# --- Align the border formatting of the newly-filled cells with their row ---
# (rows 4, 6, 8 and 10 pick up the "no-top-border" style already used by the
# neighbouring D/E columns; row 9 picks up the plain style used by B/C there)
$ws1.Range("D4:E4").Copy()
$ws1.Range("G4:H4").PasteSpecial(-4122)

$ws1.Range("D6:E6").Copy()
$ws1.Range("G6:H6").PasteSpecial(-4122)

$ws1.Range("D8:E8").Copy()
$ws1.Range("G8:H8").PasteSpecial(-4122)

$ws1.Range("D10:E10").Copy()
$ws1.Range("G10:H10").PasteSpecial(-4122)

$ws1.Range("B3:C3").Copy()
$ws1.Range("G9:H9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Restore selections to match the refreshed view ---
$ws2.Range("A2").Select()

$ws1.Activate()
$ws1.Range("H2:H11").Select()
